$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AgTests (F) and AgPosit (G) values for existing rows (revised Ag-test figures)
$ws.Range("F671").Value = 32631
$ws.Range("G671").Value = 616

$ws.Range("F678").Value = 33872
$ws.Range("G678").Value = 524

$ws.Range("F679").Value = 29476
$ws.Range("G679").Value = 519

$ws.Range("F681").Value = 26425
$ws.Range("G681").Value = 578

$ws.Range("F684").Value = 57253
$ws.Range("G684").Value = 1213

$ws.Range("F685").Value = 34472
$ws.Range("G685").Value = 1029

$ws.Range("F686").Value = 34430
$ws.Range("G686").Value = 1140

$ws.Range("F687").Value = 31477
$ws.Range("G687").Value = 1131

$ws.Range("F688").Value = 32033
$ws.Range("G688").Value = 1347

$ws.Range("F689").Value = 15722
$ws.Range("G689").Value = 1061

$ws.Range("F690").Value = 27701
$ws.Range("G690").Value = 1543

$ws.Range("F691").Value = 62380
$ws.Range("G691").Value = 2813

$ws.Range("F692").Value = 41598
$ws.Range("G692").Value = 2683

$ws.Range("F693").Value = 39443
$ws.Range("G693").Value = 2720

$ws.Range("F694").Value = 37475
$ws.Range("G694").Value = 2773

$ws.Range("F695").Value = 36852
$ws.Range("G695").Value = 3111

$ws.Range("F696").Value = 17676
$ws.Range("G696").Value = 2193

$ws.Range("F698").Value = 68511
$ws.Range("G698").Value = 5750

$ws.Range("F699").Value = 42843
$ws.Range("G699").Value = 4259

$ws.Range("F700").Value = 42992
$ws.Range("G700").Value = 4211

$ws.Range("F701").Value = 41242
$ws.Range("G701").Value = 3791

$ws.Range("F702").Value = 35770
$ws.Range("G702").Value = 3856

$ws.Range("F703").Value = 16691
$ws.Range("G703").Value = 2546

$ws.Range("F704").Value = 24522
$ws.Range("G704").Value = 3617

$ws.Range("F705").Value = 54161
$ws.Range("G705").Value = 6169

$ws.Range("F706").Value = 39919
$ws.Range("G706").Value = 4835

$ws.Range("F707").Value = 38023
$ws.Range("G707").Value = 4503

$ws.Range("F708").Value = 35039
$ws.Range("G708").Value = 4038

$ws.Range("F709").Value = 31378
$ws.Range("G709").Value = 3799

$ws.Range("F710").Value = 13291
$ws.Range("G710").Value = 2288

$ws.Range("F711").Value = 19785
$ws.Range("G711").Value = 3260

$ws.Range("F712").Value = 46677
$ws.Range("G712").Value = 5543

# Row 713 previously had no AgTests/AgPosit figures - now populated
$ws.Range("F713").Value = 34097
$ws.Range("G713").Value = 4214

# Append new day's row (2022-02-16, serial 44608)
$ws.Range("A714").Value = 44608
$ws.Range("B714").Value = 1306967
$ws.Range("C714").Value = 32730
$ws.Range("D714").Value = 18443
$ws.Range("E714").Value = 18179
$ws.Range("F714").Value = 20092
$ws.Range("G714").Value = 2497
